## real_estate_roi / Calculs.xlsx
## "Added inflation to rent, corrected the sensitivity plots"
##
## 1) Bench Apport: new "Inflation" assumption (B4/C4) and the rent
##    formula (column C) now compounds that inflation rate.
## 2) Bench Apport: the sensitivity-plot snapshot columns (K/L) are
##    refreshed to the newly computed rent/net values.
## 3) RP: the monthly-payment formula (column D, rows 38-54) is
##    re-entered as one contiguous range so Excel stores it as a single
##    shared formula (values are unchanged).
## 4) Selections / scroll position restored to match the saved file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Bench Apport"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Bench Apport")

# New row 4: Inflation assumption, styled like the other yellow
# highlighted percentage inputs (same as C2).
$ws1.Range("B4").Value = "Inflation"
$ws1.Range("C4").Value = 0.02
$ws1.Range("C4").NumberFormat = $ws1.Range("C2").NumberFormat
$ws1.Range("C4").Interior.Color = $ws1.Range("C2").Interior.Color
$ws1.Range("C4").Font.Color = $ws1.Range("C2").Font.Color

# Rent (column C) now grows every year with the new inflation rate.
$ws1.Range("C9").Formula = "=C8+(`$C`$3*(1+`$C`$4)^(`$A9-1))*12"
$ws1.Range("C10:C18").Formula = "=C9+(`$C`$3*(1+`$C`$4)^(`$A10-1))*12"

# Columns I:L are a plain-value snapshot (no formulas) used for the
# sensitivity chart; refresh K (loyer) and L (net) from the recomputed
# C/D columns so they stay in sync, same as before the edit.
for ($r = 10; $r -le 18; $r++) {
    $cVal = $ws1.Cells.Item($r, 3).Value()
    $dVal = $ws1.Cells.Item($r, 4).Value()
    $ws1.Cells.Item($r, 11).Value = $cVal
    $ws1.Cells.Item($r, 12).Value = $dVal
}

# ---------------------------------------------------------------------
# Sheet "RP"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("RP")

# Re-enter the (unchanged) monthly payment formula across the whole
# D38:D54 block so Excel consolidates it into a single shared formula
# instead of 17 independent identical ones.
$ws2.Range("D38:D54").Formula = "=-(`$C`$33/12+`$C`$35*`$G`$3/12+`$F`$33/12*(1+`$F`$34)^(`$B38-1)+`$F`$35*G38/12)"

# Restore the saved selection / scroll position on RP.
$ws2.Range("C25").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 3
$ws2.Range("D38").Select()

# ---------------------------------------------------------------------
# Restore the active sheet / selection on "Bench Apport" (it is the
# tab that is active when the file is saved).
# ---------------------------------------------------------------------
$ws1.Range("H15").Select()
